$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 759.5454999999999
$ws.Range("I41").Value = 837.5
$ws.Range("K41").Value = 837.5
$ws.Range("M41").Value = -397.5

$ws.Range("H98").Value = 10042.261
$ws.Range("I98").Value = 12592.883
$ws.Range("J98").Value = 2815.5
$ws.Range("K98").Value = 12592.883
$ws.Range("L98").Value = 2815.5
$ws.Range("M98").Value = -11094.883
$ws.Range("N98").Value = -5811.5

$ws.Range("H122").Value = 10042.261
$ws.Range("I122").Value = 12592.883
$ws.Range("J122").Value = 2815.5
$ws.Range("K122").Value = 37778.649
$ws.Range("L122").Value = 8446.5
$ws.Range("M122").Value = -35328.649
$ws.Range("N122").Value = -13346.5

$ws.Range("H132").Value = 2005.2963
$ws.Range("I132").Value = 1881.0416
$ws.Range("K132").Value = 5643.1248
$ws.Range("M132").Value = -3113.1248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 170.05882
$ws.Range("I5").Value = 143.6
$ws.Range("J5").Value = 207.85715
$ws.Range("K5").Value = 143.6
$ws.Range("L5").Value = 207.85715
$ws.Range("M5").Value = -31.59999999999999
$ws.Range("N5").Value = -431.85715

$ws.Range("H132").Value = 4862.4165
$ws.Range("I132").Value = 2483.3333
$ws.Range("J132").Value = 11999.667
$ws.Range("K132").Value = 7449.999899999999
$ws.Range("L132").Value = 35999.001
$ws.Range("M132").Value = -4919.999899999999
$ws.Range("N132").Value = -41059.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 170.05882
$ws.Range("I4").Value = 143.6
$ws.Range("J4").Value = 207.85715
$ws.Range("K4").Value = 143.6
$ws.Range("L4").Value = 207.85715
$ws.Range("M4").Value = -28.59999999999999
$ws.Range("N4").Value = -437.85715

$ws.Range("H105").Value = 5132.8887
$ws.Range("I105").Value = 5272.5654
$ws.Range("K105").Value = 5272.5654
$ws.Range("M105").Value = -3525.5654

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2368554.2
$ws.Range("I4").Value = 3461718.2
$ws.Range("K4").Value = 3461718.2
$ws.Range("M4").Value = -3461606.2

$ws.Range("H52").Value = 41998.5
$ws.Range("J52").Value = 41998.5
$ws.Range("L52").Value = 41998.5
$ws.Range("N52").Value = -42586.5

$ws.Range("H132").Value = 3414.1052
$ws.Range("I132").Value = 2144.1538
$ws.Range("K132").Value = 6432.4614
$ws.Range("M132").Value = -3902.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 603.9474
$ws.Range("I5").Value = 644.6429000000001
$ws.Range("J5").Value = 490
$ws.Range("K5").Value = 1933.9287
$ws.Range("L5").Value = 1470
$ws.Range("M5").Value = -1821.9287
$ws.Range("N5").Value = -1694

$ws.Range("H7").Value = 102.82609
$ws.Range("I7").Value = 103.3
$ws.Range("K7").Value = 309.9
$ws.Range("M7").Value = -197.9

$ws.Range("H55").Value = 473024.72
$ws.Range("I55").Value = 625145.9
$ws.Range("J55").Value = 67368.336
$ws.Range("K55").Value = 1875437.7
$ws.Range("L55").Value = 202105.008
$ws.Range("M55").Value = -1875260.7
$ws.Range("N55").Value = -202459.008

$ws.Range("H131").Value = 3178589.5
$ws.Range("I131").Value = 1391.375
$ws.Range("J131").Value = 5133788
$ws.Range("K131").Value = 4174.125
$ws.Range("L131").Value = 15401364
$ws.Range("M131").Value = 865.875
$ws.Range("N131").Value = -15411444

$ws.Range("H135").Value = 603.9474
$ws.Range("I135").Value = 644.6429000000001
$ws.Range("J135").Value = 490
$ws.Range("K135").Value = 5801.7861
$ws.Range("L135").Value = 4410
$ws.Range("M135").Value = -3266.7861
$ws.Range("N135").Value = -9480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 6299.4375
$ws.Range("J2").Value = 8.166667
$ws.Range("L2").Value = 8.166667
$ws.Range("N2").Value = -234.166667

$ws.Range("H63").Value = 500015000
$ws.Range("J63").Value = 500015000
$ws.Range("L63").Value = 500015000
$ws.Range("N63").Value = -500016372

$ws.Range("H66").Value = 500015000
$ws.Range("J66").Value = 500015000
$ws.Range("L66").Value = 1500045000
$ws.Range("N66").Value = -1500051864

$ws.Range("H97").Value = 1553.6522
$ws.Range("I97").Value = 1556.7
$ws.Range("K97").Value = 1556.7
$ws.Range("M97").Value = -1060.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1848.2972
$ws.Range("I22").Value = 816
$ws.Range("K22").Value = 816
$ws.Range("M22").Value = -521

$ws.Range("H27").Value = 1848.2972
$ws.Range("I27").Value = 816
$ws.Range("K27").Value = 816
$ws.Range("M27").Value = -709

$ws.Range("H40").Value = 4814.2856
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H46").Value = 2487.9092
$ws.Range("J46").Value = 2636.7
$ws.Range("L46").Value = 2636.7
$ws.Range("N46").Value = -3012.7

$ws.Range("H122").Value = 3992.4546
$ws.Range("I122").Value = 3817.95
$ws.Range("K122").Value = 11453.85
$ws.Range("M122").Value = -9003.849999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 93333.336
$ws.Range("J5").Value = 20000
$ws.Range("L5").Value = 20000
$ws.Range("N5").Value = -20224

$ws.Range("H74").Value = 24961
$ws.Range("I74").Value = 21525
$ws.Range("J74").Value = 25451.857
$ws.Range("K74").Value = 21525
$ws.Range("L74").Value = 25451.857
$ws.Range("M74").Value = -20589
$ws.Range("N74").Value = -27323.857

$ws.Range("H77").Value = 24961
$ws.Range("I77").Value = 21525
$ws.Range("J77").Value = 25451.857
$ws.Range("K77").Value = 64575
$ws.Range("L77").Value = 76355.571
$ws.Range("M77").Value = -59895
$ws.Range("N77").Value = -85715.571

$ws.Range("H96").Value = 5033.1665
$ws.Range("I96").Value = 3879.6
$ws.Range("K96").Value = 3879.6
$ws.Range("M96").Value = -2506.6

$ws.Range("H122").Value = 3089.081
$ws.Range("I122").Value = 1656
$ws.Range("K122").Value = 4968
$ws.Range("M122").Value = -2518

$ws.Range("H132").Value = 3129.6135
$ws.Range("I132").Value = 3132.5
$ws.Range("J132").Value = 3116.625
$ws.Range("K132").Value = 9397.5
$ws.Range("L132").Value = 9349.875
$ws.Range("M132").Value = -6867.5
$ws.Range("N132").Value = -14409.875
